$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("M2").Value = 0.9979466666666666
$ws.Range("N2").Value = 2.99384
$ws.Range("O2").Value = 0.4487505885107415
$ws.Range("P2").Value = 0.4487505885107414
$ws.Range("Q2").Value = 0.42139196152
$ws.Range("R2").Value = 3.79252765368
$ws.Range("S2").Value = 0.4487505885107415
$ws.Range("T2").Value = 0.4487505885107414

# Row 3
$ws.Range("O3").Value = 0.2623432118199488
$ws.Range("P3").Value = 0.2623432118199487
$ws.Range("S3").Value = 0.2623432118199488
$ws.Range("T3").Value = 0.2623432118199487

# Row 4
$ws.Range("M4").Value = 0.4584083333333333
$ws.Range("N4").Value = 1.375225
$ws.Range("O4").Value = 0.206134271732853
$ws.Range("P4").Value = 0.2061342717328529
$ws.Range("Q4").Value = 0.193567044425
$ws.Range("R4").Value = 1.742103399825
$ws.Range("S4").Value = 0.206134271732853
$ws.Range("T4").Value = 0.2061342717328529

# Row 5
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 0.184071
$ws.Range("N5").Value = 0.5522130000000001
$ws.Range("O5").Value = 0.08277192793645689
$ws.Range("P5").Value = 0.08277192793645688
$ws.Range("Q5").Value = 0.077725636389
$ws.Range("R5").Value = 0.6995307275010001
$ws.Range("S5").Value = 0.08277192793645689
$ws.Range("T5").Value = 0.08277192793645688
